$wb = $excel.ActiveWorkbook

# This script applies a batch data refresh to the currentAveragePrice /
# LevePrice / LeveProfit columns (H:N) on each crafting-class sheet, as
# produced by the scheduled market-data runner. Some rows lose their
# trailing LeveProfitNQ/LeveProfitHQ cell (M or N) when the corresponding
# HQ/NQ price collapses to 0, and some rows gain one back.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 38.6
$ws.Range("I38").Value = 38.6
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 115.8
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 256.2
$ws.Range("N38").ClearContents()
$ws.Range("H58").Value = 2666.6667
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 2666.6667
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 8000.000100000001
$ws.Range("N58").Value = -8300.000100000001
$ws.Range("M58").ClearContents()
$ws.Range("H87").Value = 48000
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 48000
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 48000
$ws.Range("N87").Value = -50496
$ws.Range("H90").Value = 48000
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 48000
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 144000
$ws.Range("N90").Value = -156480
$ws.Range("H107").Value = 383.8611
$ws.Range("I107").Value = 329.36
$ws.Range("J107").Value = 507.72726
$ws.Range("K107").Value = 329.36
$ws.Range("L107").Value = 507.72726
$ws.Range("M107").Value = 1590.64
$ws.Range("N107").Value = -4347.72726
$ws.Range("H138").Value = 1550.4286
$ws.Range("I138").Value = 1169.04
$ws.Range("J138").Value = 1947.7084
$ws.Range("K138").Value = 3507.12
$ws.Range("L138").Value = 5843.1252
$ws.Range("M138").Value = 1632.88
$ws.Range("N138").Value = -16123.1252

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1237.5
$ws.Range("I86").Value = 1237.5
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1237.5
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -114.5
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 1237.5
$ws.Range("I89").Value = 1237.5
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 6187.5
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -571.5
$ws.Range("N89").ClearContents()
$ws.Range("H134").Value = 5262.884
$ws.Range("I134").Value = 2381.4285
$ws.Range("J134").Value = 8013.364
$ws.Range("K134").Value = 7144.2855
$ws.Range("L134").Value = 24040.092
$ws.Range("M134").Value = -4609.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 16657
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 16657
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 16657
$ws.Range("N43").Value = -17025
$ws.Range("H101").Value = 16657
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 16657
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 16657
$ws.Range("N101").Value = -23147

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 3277.9443
$ws.Range("I54").Value = 1468
$ws.Range("J54").Value = 3639.9333
$ws.Range("K54").Value = 4404
$ws.Range("L54").Value = 10919.7999
$ws.Range("M54").Value = -3845
$ws.Range("N54").Value = -12037.7999
$ws.Range("H55").Value = 809.9
$ws.Range("I55").Value = 456
$ws.Range("J55").Value = 1045.8334
$ws.Range("K55").Value = 1368
$ws.Range("L55").Value = 3137.5002
$ws.Range("M55").Value = -1191
$ws.Range("N55").Value = -3491.5002
$ws.Range("H61").Value = 319.33334
$ws.Range("I61").Value = 304
$ws.Range("J61").Value = 350
$ws.Range("K61").Value = 912
$ws.Range("L61").Value = 1050
$ws.Range("M61").Value = -697
$ws.Range("N61").Value = -1480
$ws.Range("H64").Value = 4003
$ws.Range("I64").Value = 2004
$ws.Range("J64").Value = 10000
$ws.Range("K64").Value = 6012
$ws.Range("L64").Value = 30000
$ws.Range("M64").Value = -5742
$ws.Range("N64").Value = -30540
$ws.Range("H67").Value = 4003
$ws.Range("I67").Value = 2004
$ws.Range("J67").Value = 10000
$ws.Range("K67").Value = 6012
$ws.Range("L67").Value = 30000
$ws.Range("M67").Value = -5076
$ws.Range("N67").Value = -31872
$ws.Range("H68").Value = 974
$ws.Range("I68").Value = 300
$ws.Range("J68").Value = 1198.6666
$ws.Range("K68").Value = 900
$ws.Range("L68").Value = 3595.9998
$ws.Range("M68").Value = -89
$ws.Range("N68").Value = -5217.9998
$ws.Range("H71").Value = 974
$ws.Range("I71").Value = 300
$ws.Range("J71").Value = 1198.6666
$ws.Range("K71").Value = 2700
$ws.Range("L71").Value = 10787.9994
$ws.Range("M71").Value = 1356
$ws.Range("N71").Value = -18899.9994
$ws.Range("H74").Value = 1500
$ws.Range("I74").Value = 1500
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 4500
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -3439
$ws.Range("N74").ClearContents()
$ws.Range("H76").Value = 3071.4285
$ws.Range("I76").Value = 2000
$ws.Range("J76").Value = 3500
$ws.Range("K76").Value = 6000
$ws.Range("L76").Value = 10500
$ws.Range("M76").Value = -5617
$ws.Range("N76").Value = -11266
$ws.Range("H77").Value = 1500
$ws.Range("I77").Value = 1500
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 13500
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -8196
$ws.Range("N77").ClearContents()
$ws.Range("H79").Value = 3071.4285
$ws.Range("I79").Value = 2000
$ws.Range("J79").Value = 3500
$ws.Range("K79").Value = 6000
$ws.Range("L79").Value = 10500
$ws.Range("M79").Value = -4674
$ws.Range("N79").Value = -13152
$ws.Range("H80").Value = 1597.5
$ws.Range("I80").Value = 1200
$ws.Range("J80").Value = 1995
$ws.Range("K80").Value = 3600
$ws.Range("L80").Value = 5985
$ws.Range("M80").Value = -2664
$ws.Range("N80").Value = -7857
$ws.Range("H82").Value = 8333.666999999999
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 8333.666999999999
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 25001.001
$ws.Range("N82").Value = -25813.001
$ws.Range("H83").Value = 1597.5
$ws.Range("I83").Value = 1200
$ws.Range("J83").Value = 1995
$ws.Range("K83").Value = 10800
$ws.Range("L83").Value = 17955
$ws.Range("M83").Value = -6120
$ws.Range("N83").Value = -27315
$ws.Range("H85").Value = 8333.666999999999
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 8333.666999999999
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 25001.001
$ws.Range("N85").Value = -27809.001
$ws.Range("H86").Value = 400.4
$ws.Range("I86").Value = 400.4
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1201.2
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -15.19999999999982
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 400.4
$ws.Range("I89").Value = 400.4
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 3603.6
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 2324.4
$ws.Range("N89").ClearContents()
$ws.Range("H103").Value = 3087.75
$ws.Range("I103").Value = 274.8
$ws.Range("J103").Value = 7776
$ws.Range("K103").Value = 824.4000000000001
$ws.Range("L103").Value = 23328
$ws.Range("M103").Value = 54.59999999999991
$ws.Range("N103").Value = -25086
$ws.Range("H131").Value = 1239.4572
$ws.Range("I131").Value = 1500.6666
$ws.Range("J131").Value = 1149.0385
$ws.Range("K131").Value = 4501.9998
$ws.Range("L131").Value = 3447.1155
$ws.Range("M131").Value = 538.0002000000004
$ws.Range("N131").Value = -13527.1155

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 24088.75
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 24088.75
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 24088.75
$ws.Range("N100").Value = -26252.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4467.8184
$ws.Range("I16").Value = 2192.5715
$ws.Range("J16").Value = 8449.5
$ws.Range("K16").Value = 2192.5715
$ws.Range("L16").Value = 8449.5
$ws.Range("M16").Value = -2022.5715
$ws.Range("N16").Value = -8789.5
$ws.Range("H136").Value = 2519.1904
$ws.Range("I136").Value = 2054.9092
$ws.Range("J136").Value = 3029.9
$ws.Range("K136").Value = 6164.7276
$ws.Range("L136").Value = 9089.700000000001
$ws.Range("M136").Value = -3614.7276
$ws.Range("N136").Value = -14189.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 17010
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 17010
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 17010
$ws.Range("N92").Value = -22002
$ws.Range("H122").Value = 3637.8572
$ws.Range("I122").Value = 686.6667
$ws.Range("J122").Value = 5851.25
$ws.Range("K122").Value = 2060.0001
$ws.Range("L122").Value = 17553.75
$ws.Range("M122").Value = 389.9998999999998
$ws.Range("N122").Value = -22453.75
$ws.Range("H126").Value = 1710.5122
$ws.Range("I126").Value = 1692.0571
$ws.Range("J126").Value = 1818.1666
$ws.Range("K126").Value = 5076.1713
$ws.Range("L126").Value = 5454.4998
$ws.Range("M126").Value = -2606.1713
$ws.Range("N126").Value = -10394.4998
$ws.Range("H136").Value = 2723.375
$ws.Range("I136").Value = 3494.7144
$ws.Range("J136").Value = 2405.7646
$ws.Range("K136").Value = 10484.1432
$ws.Range("L136").Value = 7217.293799999999
$ws.Range("M136").Value = -7934.143199999999
$ws.Range("N136").Value = -12317.2938
